$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 295 (shifts old rows 295-306 down to 296-307)
$ws.Rows.Item(295).Insert()

# Populate the newly inserted row 295 with the Cleveland colo entry
$ws.Cells.Item(295, 1).Value = "CLE"
$ws.Cells.Item(295, 2).Value = "Cleveland, United States"
$ws.Cells.Item(295, 3).Value = 41.50069
$ws.Cells.Item(295, 4).Value = -81.68411999999999
$ws.Cells.Item(295, 5).Value = "US"
$ws.Cells.Item(295, 6).Value = "North America"
$ws.Cells.Item(295, 7).Value = "Cleveland"

# Match formatting of column A used throughout the data rows (bold, bordered, centered style)
$ws.Cells.Item(295, 1).Font.Bold = $true
$ws.Cells.Item(295, 1).HorizontalAlignment = -4108
$ws.Cells.Item(295, 1).VerticalAlignment = -4160
$ws.Cells.Item(295, 1).Borders.Item(1).LineStyle = 1
$ws.Cells.Item(295, 1).Borders.Item(2).LineStyle = 1
$ws.Cells.Item(295, 1).Borders.Item(3).LineStyle = 1
$ws.Cells.Item(295, 1).Borders.Item(4).LineStyle = 1
